# Update the "想去人数" (F column) counts that were refreshed by the
# gh-pages data-generation script (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) gets most of the refreshed counts.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 382
$wsExhibit.Range("F4").Value = 1606
$wsExhibit.Range("F6").Value = 23
$wsExhibit.Range("F7").Value = 411
$wsExhibit.Range("F8").Value = 142

# Sheet "全部类型" (index 4) mirrors "展览" but only a subset of rows
# were refreshed in this pass.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1606
$wsAll.Range("F8").Value = 142
$wsAll.Range("F10").Value = 475
